$d = $word.ActiveDocument

# --- "Primary" bullet used to be split across two runs ("P" + "rimary"); merge
# --- it back into a single run by doing a self Find/Replace over the visible text.
$d.Content.Find.Execute("Primary", $true, $false, $false, $false, $false, $true, 1, $false, "Primary", 2) | Out-Null

# --- Append the USERS / APPOINTMENT / PrimaryAccount DDL script as new paragraphs
# --- at the very end of the document (after the existing trailing blank lines).
$p = $d.Paragraphs.Add()
$p.Range.Text = "drop table USERS;"
$p = $d.Paragraphs.Add()
$p.Range.Text = "Create table USERS"
$p = $d.Paragraphs.Add()
$p.Range.Text = "("
$p = $d.Paragraphs.Add()
$p.Range.Text = " USER_ID NUMBER (10),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " USERNAME VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " PASSWORD VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " FIRSTNAME VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " LASTNAME VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " EMAIL VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " PHONE VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = "CONSTRAINT USER_ID_pk  PRIMARY KEY (USER_ID)"
$p = $d.Paragraphs.Add()
$p.Range.Text = ");"
$d.Paragraphs.Add() | Out-Null
$d.Paragraphs.Add() | Out-Null
$p = $d.Paragraphs.Add()
$p.Range.Text = "drop table APPOINTMENT;"
$d.Paragraphs.Add() | Out-Null
$p = $d.Paragraphs.Add()
$p.Range.Text = "Create table APPOINTMENT"
$p = $d.Paragraphs.Add()
$p.Range.Text = "("
$p = $d.Paragraphs.Add()
$p.Range.Text = "APPOINTMENT_ID NUMBER (10),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " APPOINTMENT_DATE  TIMESTAMP (9),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " APPOINTMENT_LOCATION VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " APPOINTMENT_DESCRIPTION VARCHAR2 (100),"
$p = $d.Paragraphs.Add()
$p.Range.Text = " APPOINTMENT_CONFIRMED VARCHAR (1)"
$p = $d.Paragraphs.Add()
$p.Range.Text = ");"
$d.Paragraphs.Add() | Out-Null
$p = $d.Paragraphs.Add()
$p.Range.Text = "drop table PrimaryAccount;"
$d.Paragraphs.Add() | Out-Null
$p = $d.Paragraphs.Add()
$p.Range.Text = "Create table  PrimaryAccount"
$p = $d.Paragraphs.Add()
$p.Range.Text = "("
$p = $d.Paragraphs.Add()
$p.Range.Text = "Primary_Id NUMBER (10),"
$p = $d.Paragraphs.Add()
$p.Range.Text = "Primary_account NUMBER (10),"
$p = $d.Paragraphs.Add()
$p.Range.Text = "Primary_Balance NUMBER (10,2)"
$p = $d.Paragraphs.Add()
$p.Range.Text = ");"

$d.Save()
